# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48

# Row 4
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67

# Row 7
$ws.Range("G7").Value = 3.8
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 1.85
$ws.Range("J7").Value = 4.33
$ws.Range("K7").Value = 2.25
$ws.Range("L7").Value = 2.5
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 2.05
$ws.Range("S7").Value = 1.36
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 2.1
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 13
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 9.5
$ws.Range("AK7").Value = 15
$ws.Range("AM7").Value = 151
$ws.Range("AP7").Value = 26
$ws.Range("AT7").Value = 3
$ws.Range("AU7").Value = 7.5
$ws.Range("AX7").Value = 10
$ws.Range("AZ7").Value = 34

# Row 8
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 1.98

# Row 10
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7

# Row 14
$ws.Range("G14").Value = 5.25
$ws.Range("I14").Value = 1.55
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 9.5
$ws.Range("Z14").Value = 51
$ws.Range("AC14").Value = 9.5
$ws.Range("AW14").Value = 3.5
$ws.Range("AX14").Value = 8.5

# Row 15
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 2.05
$ws.Range("K15").Value = 1.95
$ws.Range("L15").Value = 2.88
$ws.Range("O15").Value = 1.44
$ws.Range("P15").Value = 2.63
$ws.Range("S15").Value = 1.53
$ws.Range("T15").Value = 2.38
$ws.Range("U15").Value = 2.05
$ws.Range("V15").Value = 1.7
$ws.Range("W15").Value = 8.5
$ws.Range("Y15").Value = 13
$ws.Range("AA15").Value = 34
$ws.Range("AB15").Value = 41
$ws.Range("AC15").Value = 7
$ws.Range("AF15").Value = 67
$ws.Range("AG15").Value = 6
$ws.Range("AI15").Value = 9.5
$ws.Range("AL15").Value = 34
$ws.Range("AP15").Value = 34
$ws.Range("AS15").Value = 351
$ws.Range("AT15").Value = 2.38
$ws.Range("AU15").Value = 9
$ws.Range("AV15").Value = 67
$ws.Range("AX15").Value = 12
$ws.Range("AY15").Value = 26
$ws.Range("BB15").Value = 251
$ws.Range("BD15").Value = 126

# Row 16
$ws.Range("Q16").Value = 1.93
$ws.Range("R16").Value = 1.93

# Row 18
$ws.Range("G18").Value = 7.4
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 1.35
$ws.Range("J18").Value = 6.1
$ws.Range("K18").Value = 2.67
$ws.Range("S18").Value = 1.22
$ws.Range("T18").Value = 3.85
$ws.Range("U18").Value = 1.55
$ws.Range("V18").Value = 2.3
$ws.Range("AD18").Value = 11
$ws.Range("AH18").Value = 9.25
$ws.Range("AI18").Value = 8.5
$ws.Range("AJ18").Value = 10.5
$ws.Range("AL18").Value = 17.5
$ws.Range("AP18").Value = 29
$ws.Range("AR18").Value = 175
$ws.Range("AT18").Value = 3.85
$ws.Range("AU18").Value = 7.1
$ws.Range("AX18").Value = 6.1
$ws.Range("AY18").Value = 11.75
$ws.Range("AZ18").Value = 14.5

$wb.Save()